# Generate Report for Handoff
# Updates the localization-status workbook with a freshly generated
# report: new source-file GUID, new xliff content hashes, and new
# handoff/generation timestamps.

$wb = $excel.ActiveWorkbook

# old GUID being replaced: a4836735-19cc-4df8-aa28-2163973ead67
$newGuid = "cc54d45f-e499-49a1-8ed4-e1bdfd642a56"
$newHash = "feba180e7bd4ad3af8137b31740a74d5d3a92386"

# --- "Overview" sheet ---------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("A2").Value = "$newGuid.md"
$wsOverview.Range("B2").Value = "e2e\$newGuid.md"
$wsOverview.Range("G2").Value = "2016-08-23 00:54:39"
foreach ($h in $wsOverview.Hyperlinks) {
    $h.TextToDisplay = "e2e\$newGuid.md"
}

# --- "zh-cn" sheet --------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("A2").Value = "$newGuid.md"
$wsZhCn.Range("G2").Value = "$newGuid.$newHash.zh-cn.xlf"
$wsZhCn.Range("H2").Value = "2016-08-23 00:54:34"
foreach ($h in $wsZhCn.Hyperlinks) {
    $h.TextToDisplay = "$newGuid.md"
}

# --- "de-de" sheet --------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("A2").Value = "$newGuid.md"
$wsDeDe.Range("G2").Value = "$newGuid.$newHash.de-de.xlf"
$wsDeDe.Range("H2").Value = "2016-08-23 00:54:39"
foreach ($h in $wsDeDe.Hyperlinks) {
    $h.TextToDisplay = "$newGuid.md"
}
